# Swap the "Bottle ID" (column B) and "Glass ID" (column C) columns on
# Sheet1: every value that lived in B now lives in C and vice versa,
# headers included, and the column widths travel with the data.
#
# Cutting column C and inserting it before column B performs an exact
# swap (values, shared-string refs, and stored column widths) in one
# shot, which is both simpler and more faithful than rewriting every
# cell by hand.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns("C").Cut() | Out-Null
$ws.Columns("B").Insert() | Out-Null

# Update the active selection to reflect where the user ended up after
# the edit (also drops any stale scrolled-down "topLeftCell" position).
$ws.Range("E19").Select() | Out-Null
